# Auto-generated update: recompute "剩余" (remaining days) and roll over
# "开始时间" (start date) for each row, as of the new reference date
# 2025-11-07 (commit run 2025-11-06 23:12:06, rolling into the next day).
#
# Business rule recovered from the diff: end date = F (start, YYYYMMDD) + D
# (total days). If that end date has already passed (<= today), the cycle
# resets: F becomes today and E (remaining) becomes D again. Otherwise E is
# simply the number of days left until the end date; F is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2; E=10; F=20251103},
    @{Row=3; E=10; F=20251103},
    @{Row=4; E=10; F=20251103},
    @{Row=5; E=8; F=20251105},
    @{Row=6; E=10; F=20251103},
    @{Row=7; E=8; F=20251105},
    @{Row=8; E=10; F=20251103},
    @{Row=9; E=8; F=20251105},
    @{Row=10; E=3; F=20251103},
    @{Row=11; E=10; F=20251103},
    @{Row=12; E=8; F=20251105},
    @{Row=13; E=10; F=20251103},
    @{Row=14; E=10; F=20251103},
    @{Row=15; E=10; F=20251103},
    @{Row=16; E=2; F=20251030},
    @{Row=17; E=8; F=20251105},
    @{Row=18; E=1; F=20251029},
    @{Row=19; E=1; F=20251029},
    @{Row=20; E=1; F=20251029},
    @{Row=21; E=1; F=20251029},
    @{Row=22; E=8; F=20251105},
    @{Row=23; E=8; F=20251105},
    @{Row=24; E=8; F=20251105},
    @{Row=25; E=8; F=20251105},
    @{Row=26; E=8; F=20251105},
    @{Row=27; E=4; F=20251104},
    @{Row=28; E=1; F=20251029},
    @{Row=29; E=1; F=20251029},
    @{Row=30; E=1; F=20251029},
    @{Row=31; E=1; F=20251029},
    @{Row=32; E=1; F=20251029},
    @{Row=33; E=1; F=20251029},
    @{Row=34; E=1; F=20251029},
    @{Row=35; E=1; F=20251029},
    @{Row=37; E=1; F=20251029},
    @{Row=38; E=1; F=20251029},
    @{Row=39; E=1; F=20251029},
    @{Row=40; E=3; F=20251103},
    @{Row=41; E=3; F=20251103},
    @{Row=42; E=1; F=20251029},
    @{Row=43; E=8; F=20251105},
    @{Row=44; E=3; F=20251103},
    @{Row=45; E=8; F=20251105},
    @{Row=46; E=3; F=20251103},
    @{Row=47; E=1; F=20251029},
    @{Row=48; E=3; F=20251103},
    @{Row=49; E=4; F=20251104},
    @{Row=50; E=6; F=20251103},
    @{Row=51; E=6; F=20251103},
    @{Row=52; E=6; F=20251103},
    @{Row=53; E=6; F=20251103},
    @{Row=54; E=6; F=20251103},
    @{Row=55; E=6; F=20251103},
    @{Row=56; E=6; F=20251103},
    @{Row=57; E=6; F=20251103},
    @{Row=58; E=10; F=20251107},
    @{Row=59; E=10; F=20251107},
    @{Row=60; E=10; F=20251107},
    @{Row=61; E=4; F=20251104},
    @{Row=62; E=10; F=20251107},
    @{Row=63; E=10; F=20251107},
    @{Row=64; E=10; F=20251107},
    @{Row=65; E=1; F=20251029},
    @{Row=66; E=1; F=20251029},
    @{Row=67; E=1; F=20251029},
    @{Row=68; E=1; F=20251029},
    @{Row=69; E=1; F=20251029},
    @{Row=70; E=2; F=20251030},
    @{Row=71; E=2; F=20251030},
    @{Row=72; E=2; F=20251030},
    @{Row=73; E=2; F=20251030},
    @{Row=74; E=2; F=20251030},
    @{Row=75; E=2; F=20251030},
    @{Row=76; E=2; F=20251030},
    @{Row=77; E=5; F=20251102},
    @{Row=78; E=5; F=20251102},
    @{Row=79; E=5; F=20251102},
    @{Row=80; E=5; F=20251102},
    @{Row=81; E=5; F=20251102},
    @{Row=82; E=5; F=20251102},
    @{Row=83; E=5; F=20251102},
    @{Row=84; E=5; F=20251102},
    @{Row=85; E=5; F=20251102},
    @{Row=86; E=5; F=20251102},
    @{Row=87; E=3; F=20251103},
    @{Row=88; E=3; F=20251103},
    @{Row=89; E=3; F=20251103},
    @{Row=90; E=3; F=20251103},
    @{Row=91; E=8; F=20251105},
    @{Row=92; E=3; F=20251103},
    @{Row=93; E=5; F=20251102},
    @{Row=94; E=6; F=20251106},
    @{Row=95; E=4; F=20251101},
    @{Row=96; E=2; F=20251030},
    @{Row=97; E=2; F=20251030},
    @{Row=98; E=2; F=20251030},
    @{Row=99; E=2; F=20251030}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 5).Value = $u.E
    $ws.Cells.Item($u.Row, 6).Value = $u.F
}
